$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# GETs sheet: add 8 new test-case rows (6-13) for the "debin2" and
# "debin4" methods (mirroring the existing "debin5" rows 3-5).
# ------------------------------------------------------------------
$ws = $wb.Worksheets.Item("GETs")
$ws.Activate()

$zwsp = [char]0x200b

# Row 6
$ws.Range("A6").Value = 59716
$ws.Range("B6").Value = "QRDebin*->ConfirmaDebito*->debindebin2*"
$ws.Range("C6").Value = '"operacion":{"detalle":{"importe":1000}}}|"operacion":{"comprador":{"cuit":"23000009989","cuenta":{"cbu":"9985340400000000000529"}},"detalle":{"ori_terminal":"","ori_adicional":"","moneda":"032","importe":1000}}}|"id":"debin.id"'
$ws.Range("C6").WrapText = $true
$ws.Range("E6").Value = "GET /apiDebinV1" + $zwsp + "/Debin" + $zwsp + "/Debin2" + $zwsp + "/{id} - Nuevos campos en JSON"

# Row 7
$ws.Range("A7").Value = 59718
$ws.Range("B7").Value = "debindebin2*"
$ws.Range("C7").Value = '"id":"M67REZ8NP1E680124KVGOP"'
$ws.Range("C7").WrapText = $true
$ws.Range("E7").Value = "GET /apiDebinV1" + $zwsp + "/Debin" + $zwsp + "/Debin2" + $zwsp + "/{id} - Devolucion Parcial"

# Row 8
$ws.Range("A8").Value = 59720
$ws.Range("B8").Value = "debindebin2*"
$ws.Range("C8").Value = '"id":"M67REZ8NP1E64ZG24KVMM"'
$ws.Range("D8").Value = '{"StatusCode":200,"Mensaje":{"respuesta": {"codigo":"83","descripcion":"DEBIN INEXISTENTE"}}}'
$ws.Range("D8").WrapText = $true
$ws.Range("D8").HorizontalAlignment = -4131
$ws.Range("D8").VerticalAlignment = -4160
$ws.Range("E8").Value = "GET /apiDebinV1" + $zwsp + "/Debin" + $zwsp + "/Debin2" + $zwsp + "/{id} - ID HASH inexistente"

# Row 9
$ws.Range("A9").Value = 53394
$ws.Range("B9").Value = "QRDebin*->debindebin2*"
$ws.Range("C9").Value = '"operacion":{"detalle":{"importe":1000}}}|"id":"debin.id"'
$ws.Range("C9").WrapText = $true
$ws.Range("E9").Value = "Get Debin 2 --> consultar por un Debin QR"

# Row 10
$ws.Range("A10").Value = 53395
$ws.Range("B10").Value = "CashOut*->debindebin2*"
$ws.Range("C10").Value = '"objeto":{"tipo":"CASHOUT"},"credito":{"cuit":"20956746117","banco":"000","sucursal":"0213","cuenta":{"cbu":"0000213699900070000000"},"titular":""},"debito":{"cuit":"20333048494","banco":"998","sucursal":"8851","cuenta":{"cbu":"9988851800000000000628"},"titular":"PRUEBAS COELSA CASHOUT"},"importe":{"importe":10}}|"id":"debin.id"'
$ws.Range("C10").WrapText = $true
$ws.Range("C10").HorizontalAlignment = -4131
$ws.Range("C10").VerticalAlignment = -4108
$ws.Range("E10").Value = "Get Debin 2 --> consultar por un Cashout"

# Row 11
$ws.Range("A11").Value = 59542
$ws.Range("B11").Value = "QRDebin*->ConfirmaDebito*->debindebin4*"
$ws.Range("C11").Value = '"operacion":{"detalle":{"importe":1000}}}|"operacion":{"comprador":{"cuit":"23000009989","cuenta":{"cbu":"9985340400000000000529"}},"detalle":{"ori_terminal":"","ori_adicional":"","moneda":"032","importe":1000}}}|"id":"debin.id"'
$ws.Range("C11").WrapText = $true
$ws.Range("E11").Value = "GET /apiDebinV1" + $zwsp + "/Debin" + $zwsp + "/Debin5" + $zwsp + "/{id} - Nuevos campos en JSON"

# Row 12
$ws.Range("A12").Value = 59545
$ws.Range("B12").Value = "debindebin4*"
$ws.Range("C12").Value = '"id":"M67REZ8NP1E680124KVGOP"'
$ws.Range("C12").WrapText = $true
$ws.Range("E12").Value = "GET /apiDebinV1" + $zwsp + "/Debin" + $zwsp + "/Debin5" + $zwsp + "/{id} - Devolucion Parcial"

# Row 13
$ws.Range("A13").Value = 59551
$ws.Range("B13").Value = "debindebin4*"
$ws.Range("C13").Value = '"id":"M67REZ8NP1E64ZG24KVMM"'
$ws.Range("D13").Value = '{"StatusCode":200,"Mensaje":{"respuesta": {"codigo":"83","descripcion":"DEBIN INEXISTENTE"}}}'
$ws.Range("D13").WrapText = $true
$ws.Range("D13").HorizontalAlignment = -4131
$ws.Range("D13").VerticalAlignment = -4160
$ws.Range("E13").Value = "GET /apiDebinV1" + $zwsp + "/Debin" + $zwsp + "/Debin5" + $zwsp + "/{id} - ID HASH inexistente"

$ws.Rows.Item(6).RowHeight = 45
$ws.Rows.Item(8).RowHeight = 30
$ws.Rows.Item(10).RowHeight = 60
$ws.Rows.Item(11).RowHeight = 45
$ws.Rows.Item(13).RowHeight = 30

# ------------------------------------------------------------------
# CONTRACARGOQR sheet: selection moved from B4 to B7.
# ------------------------------------------------------------------
$wsContra = $wb.Worksheets.Item("CONTRACARGOQR")
$wsContra.Activate()
$wsContra.Range("B7").Select()

# ------------------------------------------------------------------
# CASHOUT sheet: selection moved from A4 to B4:C4.
# ------------------------------------------------------------------
$wsCashout = $wb.Worksheets.Item("CASHOUT")
$wsCashout.Activate()
$wsCashout.Range("B4:C4").Select()

# ------------------------------------------------------------------
# Leave GETs as the active/selected sheet, with A13 selected.
# ------------------------------------------------------------------
$ws.Activate()
$ws.Range("A13").Select()
